# Changes after drive practice 3/9
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the measured velocity value for row 3 (C3: 9750 -> 10200)
$ws.Range("C3").Value = 10200

# Keep the chart's plotted series in sync with the updated source data.
try {
    $chartObj = $ws.ChartObjects().Item(1)
    $chart = $chartObj.Chart
    $chart.SeriesCollection(1).Values = $ws.Range("C2:C9")
} catch {
}

$excel.Calculate()

# Move the active selection to C4 (previously A3:C3)
$ws.Range("C4").Select()
